$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-09 04:28:55"
$wsZhCn.Range("G3").Value = "2016-01-09 04:29:38"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-09 04:29:04"
$wsDeDe.Range("G3").Value = "2016-01-09 04:29:55"
